# Update NATMI Efna5-Ephb6 LR-pairs results with newly computed TPM values.
# The sender/target clusters are now limited to ECs, FAPs, MuSCs (the
# "Resolving-Mac" cluster combinations are removed), so the data block
# shrinks from 12 rows (3x4 combinations) to 9 rows (3x3 combinations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three rows that referenced the now-removed "Resolving-Mac"
# target cluster (old rows 11-13); remaining rows shift up automatically.
$ws.Range("A11:T13").EntireRow.Delete()

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Ephb6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1728506666666667
$ws.Cells.Item(2, 8).Value = 0.518552
$ws.Cells.Item(2, 9).Value = 0.0840503369699626
$ws.Cells.Item(2, 10).Value = 0.0840503369699626
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.168144
$ws.Cells.Item(2, 14).Value = 0.504432
$ws.Cells.Item(2, 15).Value = 0.05446245276675245
$ws.Cells.Item(2, 16).Value = 0.05446245276675245
$ws.Cells.Item(2, 17).Value = 0.029063802496
$ws.Cells.Item(2, 18).Value = 0.261574222464
$ws.Cells.Item(2, 19).Value = 0.004577587507256215
$ws.Cells.Item(2, 20).Value = 0.004577587507256215
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Ephb6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1728506666666667
$ws.Cells.Item(3, 8).Value = 0.518552
$ws.Cells.Item(3, 9).Value = 0.0840503369699626
$ws.Cells.Item(3, 10).Value = 0.0840503369699626
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.279928333333333
$ws.Cells.Item(3, 14).Value = 3.839785
$ws.Cells.Item(3, 15).Value = 0.4145734394268892
$ws.Cells.Item(3, 16).Value = 0.4145734394268892
$ws.Cells.Item(3, 17).Value = 0.2212364657022222
$ws.Cells.Item(3, 18).Value = 1.99112819132
$ws.Cells.Item(3, 19).Value = 0.03484503728262642
$ws.Cells.Item(3, 20).Value = 0.03484503728262641
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Ephb6"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1728506666666667
$ws.Cells.Item(4, 8).Value = 0.518552
$ws.Cells.Item(4, 9).Value = 0.0840503369699626
$ws.Cells.Item(4, 10).Value = 0.0840503369699626
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.639265666666667
$ws.Cells.Item(4, 14).Value = 4.917797
$ws.Cells.Item(4, 15).Value = 0.5309641078063584
$ws.Cells.Item(4, 16).Value = 0.5309641078063583
$ws.Cells.Item(4, 17).Value = 0.2833481633271112
$ws.Cells.Item(4, 18).Value = 2.550133469944
$ws.Cells.Item(4, 19).Value = 0.04462771218007997
$ws.Cells.Item(4, 20).Value = 0.04462771218007996
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Ephb6"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.367717666666667
$ws.Cells.Item(5, 8).Value = 4.103153
$ws.Cells.Item(5, 9).Value = 0.6650661694281633
$ws.Cells.Item(5, 10).Value = 0.6650661694281633
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.168144
$ws.Cells.Item(5, 14).Value = 0.504432
$ws.Cells.Item(5, 15).Value = 0.05446245276675245
$ws.Cells.Item(5, 16).Value = 0.05446245276675245
$ws.Cells.Item(5, 17).Value = 0.229973519344
$ws.Cells.Item(5, 18).Value = 2.069761674096
$ws.Cells.Item(5, 19).Value = 0.03622113483924633
$ws.Cells.Item(5, 20).Value = 0.03622113483924633
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Ephb6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.367717666666667
$ws.Cells.Item(6, 8).Value = 4.103153
$ws.Cells.Item(6, 9).Value = 0.6650661694281633
$ws.Cells.Item(6, 10).Value = 0.6650661694281633
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.279928333333333
$ws.Cells.Item(6, 14).Value = 3.839785
$ws.Cells.Item(6, 15).Value = 0.4145734394268892
$ws.Cells.Item(6, 16).Value = 0.4145734394268892
$ws.Cells.Item(6, 17).Value = 1.750580593567222
$ws.Cells.Item(6, 18).Value = 15.755225342105
$ws.Cells.Item(6, 19).Value = 0.2757187693062999
$ws.Cells.Item(6, 20).Value = 0.2757187693062999
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Ephb6"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.367717666666667
$ws.Cells.Item(7, 8).Value = 4.103153
$ws.Cells.Item(7, 9).Value = 0.6650661694281633
$ws.Cells.Item(7, 10).Value = 0.6650661694281633
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.639265666666667
$ws.Cells.Item(7, 14).Value = 4.917797
$ws.Cells.Item(7, 15).Value = 0.5309641078063584
$ws.Cells.Item(7, 16).Value = 0.5309641078063583
$ws.Cells.Item(7, 17).Value = 2.242052612660111
$ws.Cells.Item(7, 18).Value = 20.178473513941
$ws.Cells.Item(7, 19).Value = 0.3531262652826171
$ws.Cells.Item(7, 20).Value = 0.353126265282617
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Efna5"
$ws.Cells.Item(8, 3).Value = "Ephb6"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5159453333333334
$ws.Cells.Item(8, 8).Value = 1.547836
$ws.Cells.Item(8, 9).Value = 0.2508834936018741
$ws.Cells.Item(8, 10).Value = 0.2508834936018741
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.168144
$ws.Cells.Item(8, 14).Value = 0.504432
$ws.Cells.Item(8, 15).Value = 0.05446245276675245
$ws.Cells.Item(8, 16).Value = 0.05446245276675245
$ws.Cells.Item(8, 17).Value = 0.08675311212799999
$ws.Cells.Item(8, 18).Value = 0.780778009152
$ws.Cells.Item(8, 19).Value = 0.01366373042024991
$ws.Cells.Item(8, 20).Value = 0.01366373042024991
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Efna5"
$ws.Cells.Item(9, 3).Value = "Ephb6"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5159453333333334
$ws.Cells.Item(9, 8).Value = 1.547836
$ws.Cells.Item(9, 9).Value = 0.2508834936018741
$ws.Cells.Item(9, 10).Value = 0.2508834936018741
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.279928333333333
$ws.Cells.Item(9, 14).Value = 3.839785
$ws.Cells.Item(9, 15).Value = 0.4145734394268892
$ws.Cells.Item(9, 16).Value = 0.4145734394268892
$ws.Cells.Item(9, 17).Value = 0.6603730505844445
$ws.Cells.Item(9, 18).Value = 5.94335745526
$ws.Cells.Item(9, 19).Value = 0.1040096328379629
$ws.Cells.Item(9, 20).Value = 0.1040096328379629
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Efna5"
$ws.Cells.Item(10, 3).Value = "Ephb6"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5159453333333334
$ws.Cells.Item(10, 8).Value = 1.547836
$ws.Cells.Item(10, 9).Value = 0.2508834936018741
$ws.Cells.Item(10, 10).Value = 0.2508834936018741
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.639265666666667
$ws.Cells.Item(10, 14).Value = 4.917797
$ws.Cells.Item(10, 15).Value = 0.5309641078063584
$ws.Cells.Item(10, 16).Value = 0.5309641078063583
$ws.Cells.Item(10, 17).Value = 0.8457714708102223
$ws.Cells.Item(10, 18).Value = 7.611943237292
$ws.Cells.Item(10, 19).Value = 0.1332101303436613
$ws.Cells.Item(10, 20).Value = 0.1332101303436613
